# Add a "Total Price" column (K) to the BOM sheet: bold header matching the
# other headers, a per-row Quantity*Unit-Price formula for every part row,
# and a grand-total SUM underneath.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell K1: "Total Price", bold like the rest of row 1 -----------
$ws.Range("K1").Value2 = "Total Price"
$ws.Range("K1").Font.Bold = $true

# --- K2 = I2*B2  (Unit Price * Quantity) for the first part row -----------
$ws.Range("K2").Formula = "=I2*B2"
$ws.Range("K2").Style = "Normal"

# --- K3:K28 = I{row}*B{row} for the remaining part rows --------------------
$ws.Range("K3:K28").Formula = "=I3*B3"
$ws.Range("K3:K28").Style = "Normal"

# --- K29 = SUM(K2:K28) grand total -----------------------------------------
$ws.Range("K29").Formula = "=SUM(K2:K28)"

# --- Restore the view: scrolled to C3, with G11 selected -------------------
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("G11").Select()
